$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01253208636536152, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1, 0.7365606309410384)
    3 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367)
    4 = @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.598097515653722)
    5 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    6 = @(0.0006075818656279264, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1, 0.7246361264413048)
    7 = @(0.04172184405617529, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0, 0.7657503886318522)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
